$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.21888754210301897
$ws.Range("B1").Value = 0.21854790196423579
$ws.Range("A2").Value = -0.16010721563938191
$ws.Range("B2").Value = 0.1592078006129638
$ws.Range("A3").Value = -0.10949479055234868
$ws.Range("B3").Value = 0.10929005038771855
$ws.Range("A4").Value = -0.10129005042773187
$ws.Range("B4").Value = 0.10093811789773
$ws.Range("A5").Value = -0.09793811792061824
$ws.Range("B5").Value = 0.096752962713377855
$ws.Range("A6").Value = 0.020391242916996788
$ws.Range("B6").Value = -0.020408263044826569
$ws.Range("A7").Value = 0.03040826298832755
$ws.Range("B7").Value = -0.030427499125466717
$ws.Range("A8").Value = 0.040427499070240902
$ws.Range("B8").Value = -0.040578099862036865
$ws.Range("A9").Value = 0.042578099840962835
$ws.Range("B9").Value = -0.042760727704305612
$ws.Range("A10").Value = 0.044760727686130153
$ws.Range("B10").Value = -0.044778949085250375
$ws.Range("A11").Value = 0.047778949063494558
$ws.Range("B11").Value = -0.047820984054798288
$ws.Range("A12").Value = 0.0015796002979264756
$ws.Range("B12").Value = -0.0015872562660237399
$ws.Range("A13").Value = 0.0050872562435131385
$ws.Range("B13").Value = -0.0050888747300659176
$ws.Range("A14").Value = 0.013088874690034835
$ws.Range("B14").Value = -0.013096504527539743
$ws.Range("A15").Value = -0.0080501448628291783
$ws.Range("B15").Value = 0.0080326019997514209
$ws.Range("A16").Value = -0.0060326020161709515
$ws.Range("B16").Value = 0.0060033757636928087
$ws.Range("A17").Value = -0.004003375780583518
$ws.Range("B17").Value = 0.0039999999751332282
$ws.Range("A18").Value = -0.052440686349385146
$ws.Range("B18").Value = 0.052361413222051567
$ws.Range("A19").Value = -0.048361413239158324
$ws.Range("B19").Value = 0.04779846361051332
$ws.Range("A20").Value = -0.0080172130206950243
$ws.Range("B20").Value = 0.0080057731311864444
$ws.Range("A21").Value = -0.0040057731533176266
$ws.Range("B21").Value = 0.0039999999777320383
$ws.Range("A22").Value = -0.045713010082430117
$ws.Range("B22").Value = 0.045499310304441209
$ws.Range("A23").Value = -0.040499310331493454
$ws.Range("B23").Value = 0.040099122960925548
$ws.Range("A24").Value = -0.02009912304919137
$ws.Range("B24").Value = 0.019999999910611521
$ws.Range("A25").Value = -0.11514420565563732
$ws.Range("B25").Value = 0.1149588233893617
$ws.Range("A26").Value = -0.11245882341803259
$ws.Range("B26").Value = 0.1122170905832558
$ws.Range("A27").Value = -0.091966553908140902
$ws.Range("B27").Value = 0.091014955146782928
$ws.Range("A28").Value = -0.089014955183665201
$ws.Range("B28").Value = 0.08836338976826319
$ws.Range("A29").Value = -0.081363389830611865
$ws.Range("B29").Value = 0.081174361785157245
$ws.Range("A30").Value = -0.021174362058645535
$ws.Range("B30").Value = 0.021024094949874783
$ws.Range("A31").Value = -0.014024095016962335
$ws.Range("B31").Value = 0.014001218381114455
$ws.Range("A32").Value = -0.0040012184602851164
$ws.Range("B32").Value = 0.0039999999443285361
